$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controlador de projetos")

# Row 14: fill in the "Real Término" (J14) date; K14 formula recalculates automatically
$ws.Range("J14").Value = 45723

# Row 15: Proposta Analítca
$ws.Range("B15").Value = "Proposta Anaíltca"
$ws.Range("C15").Value = "Todos"
$ws.Range("D15").Value = "Segunda Entrega"
$ws.Range("E15").Value = "Documento"
$ws.Range("F15").Value = 45719
$ws.Range("G15").Value = 45741

# Row 16: Análise Exploratória
$ws.Range("B16").Value = "Análise Exploratória"
$ws.Range("C16").Value = "Todos"
$ws.Range("D16").Value = "Segunda Entrega"
$ws.Range("E16").Value = "Análise"
$ws.Range("F16").Value = 45719
$ws.Range("G16").Value = 45741

# Row 17: Scripts Análise Exploratória
$ws.Range("B17").Value = "Scripts Análise Exploratória"
$ws.Range("C17").Value = "Todos"
$ws.Range("D17").Value = "Segunda Entrega"
$ws.Range("E17").Value = "Código"
$ws.Range("F17").Value = 45719
$ws.Range("G17").Value = 45741

# Row 18: Esboço Sotory Telling
$ws.Range("B18").Value = "Esboço Sotory Telling"
$ws.Range("C18").Value = "Todos"
$ws.Range("D18").Value = "Terceira Entrega"
$ws.Range("E18").Value = "Apresentação"
$ws.Range("F18").Value = 45747
$ws.Range("G18").Value = 45772

# Row 19: Revisão Scripts A.E.
$ws.Range("B19").Value = "Revisão Scripts A.E."
$ws.Range("C19").Value = "Todos"
$ws.Range("D19").Value = "Terceira Entrega"
$ws.Range("E19").Value = "Código"
$ws.Range("F19").Value = 45747
$ws.Range("G19").Value = 45772

# Row 20: A. E. no relatório
$ws.Range("B20").Value = "A. E. no relatório"
$ws.Range("C20").Value = "Todos"
$ws.Range("D20").Value = "Terceira Entrega"
$ws.Range("E20").Value = "Documento"
$ws.Range("F20").Value = 45747
$ws.Range("G20").Value = 45772

# Row 21: Video Story Telling
$ws.Range("B21").Value = "Video Story Telling"
$ws.Range("C21").Value = "Todos"
$ws.Range("D21").Value = "Quarta Entrega"
$ws.Range("E21").Value = "Apresentação"
$ws.Range("F21").Value = 45775
$ws.Range("G21").Value = 45800

# Row 22: Relatório Final
$ws.Range("B22").Value = "Relatório Final"
$ws.Range("C22").Value = "Todos"
$ws.Range("D22").Value = "Quarta Entrega"
$ws.Range("E22").Value = "Documento"
$ws.Range("F22").Value = 45775
$ws.Range("G22").Value = 45800

# Header note showing the last update date (new shared string added last)
$ws.Range("C2").Value = "Atualizado em 2025-03-07"

# Update the active selection to C3
$ws.Range("C3").Select()
